$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update the summary header figures (fewer workers / periods after the
# database refresh, and a recalculated total "Valor Mora").
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 409679
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 9

# ---------------------------------------------------------------------------
# Remove one worker row from the detail table (data refresh removed
# EFRAIN ISAAC MARTINEZ MARTINEZ entirely). Deleting any one of the
# uniformly-styled rows 16-25 shifts the remaining rows up, carrying the
# special "last row" bottom-border formatting down onto the new last row
# automatically.
# ---------------------------------------------------------------------------
$ws.Rows(17).Delete()

# ---------------------------------------------------------------------------
# Rewrite the worker detail table (rows 16-25) with the refreshed data
# pulled from the updated database export.
# ---------------------------------------------------------------------------
$data = @(
    @("PPT", "1287309",   "WLADIMIR ALEXANDER GARCIA PERAZA", "2302", 46400, 1423500),
    @("PPT", "1287309",   "WLADIMIR ALEXANDER GARCIA PERAZA", "2303", 46400, 1423500),
    @("PPT", "1287309",   "WLADIMIR ALEXANDER GARCIA PERAZA", "2304", 46400, 1423500),
    @("CC",  "20429869",  "ALEXANDER MOISES ANILLO MONTES",   "2304", 18372, 1531000),
    @("PPT", "1287309",   "WLADIMIR ALEXANDER GARCIA PERAZA", "2305", 46400, 1423500),
    @("PPT", "1287309",   "WLADIMIR ALEXANDER GARCIA PERAZA", "2306", 46400, 1423500),
    @("PPT", "1287309",   "WLADIMIR ALEXANDER GARCIA PERAZA", "2307", 46400, 1423500),
    @("PPT", "1287309",   "WLADIMIR ALEXANDER GARCIA PERAZA", "2308", 46400, 1423500),
    @("PPT", "1287309",   "WLADIMIR ALEXANDER GARCIA PERAZA", "2309", 46400, 1423500),
    @("CC",  "1143375685","EDGARDO LUIS VERGARA CONTO",       "2311", 20107, 1160000)
)

$row = 16
foreach ($entry in $data) {
    $ws.Range("B$row").Value = $entry[0]
    $ws.Range("C$row").Value = $entry[1]
    $ws.Range("D$row").Value = $entry[2]
    $ws.Range("E$row").Value = $entry[3]
    $ws.Range("F$row").Value = $entry[4]
    $ws.Range("G$row").Value = $entry[5]
    $row = $row + 1
}
